# Dev IV Project Rubric - end of lab edit
# Fills in the missing milestone grade (Roman numeral + X mark) for two
# rows (6 and 10) that were previously left blank, and restores the view
# to the top of the sheet with the selection on E25.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 6: mark milestone III with an X (Student's 3rd milestone score)
$ws.Range("E6").Value = "III"
$ws.Range("F6").Value = "X"

# Row 10: mark milestone III with an X as well
$ws.Range("E10").Value = "III"
$ws.Range("F10").Value = "X"

# Restore the sheet view to the top-left area, with E25 selected
$excel.ActiveWindow.ScrollColumn = 3
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("E25").Select()
